# Generate Report for Handback
# Adds the handback record for file 32a2d89a-a2e6-4f6a-b723-17d9cd919a75.md
# as a new row (row 4) on the "Overview", "zh-cn" and "de-de" sheets, and
# grows the backing tables/dimensions to match.

$wb = $excel.ActiveWorkbook

$guid        = "32a2d89a-a2e6-4f6a-b723-17d9cd919a75"
$mdName      = "$guid.md"
$srcPath     = "e2e\$guid.md"
$xliffSha    = "9b742254ffe519d9b27b3f0a4a3eb0f2e6c0ede7"
$zhXlf       = "$guid.$xliffSha.zh-cn.xlf"
$deXlf       = "$guid.$xliffSha.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"

$hoDateZh  = "2016-08-30 06:47:04"
$hbDateZh  = "2016-08-30 06:47:36"
$hoDateDe  = "2016-08-30 06:47:12"
$hbDateDe  = "2016-08-30 06:47:43"
$latestDate = "2016-08-30 06:47:12"

$urlBase    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b742254ffe519d9b27b3f0a4a3eb0f2e6c0ede7/e2e/$mdName"
$urlZhCn    = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9b742254ffe519d9b27b3f0a4a3eb0f2e6c0ede7/e2e/$mdName"
$urlDeDe    = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9b742254ffe519d9b27b3f0a4a3eb0f2e6c0ede7/e2e/$mdName"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 0xED9564
}

function Style-AsDate($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Overview sheet -> row 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $srcPath
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = $latestDate

Style-AsHyperlink $wsOverview.Range("B4")
Style-AsDate $wsOverview.Range("G4")

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $urlBase, "", "", $srcPath) | Out-Null

$wsOverview.Range("A1:G4").Columns.AutoFit() | Out-Null
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4")) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet -> row 4
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $statusInSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $hoDateZh
$wsZh.Range("I4").Value = $mdName
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $hbDateZh
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

Style-AsHyperlink $wsZh.Range("A4")
Style-AsHyperlink $wsZh.Range("I4")
Style-AsDate $wsZh.Range("H4")
Style-AsDate $wsZh.Range("K4")

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $urlBase, "", "", $mdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $urlZhCn, "", "", $mdName) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4")) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet -> row 4
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $statusInSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $hoDateDe
$wsDe.Range("I4").Value = $mdName
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $hbDateDe
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

Style-AsHyperlink $wsDe.Range("A4")
Style-AsHyperlink $wsDe.Range("I4")
Style-AsDate $wsDe.Range("H4")
Style-AsDate $wsDe.Range("K4")

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $urlBase, "", "", $mdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $urlDeDe, "", "", $mdName) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4")) | Out-Null

Write-Host "Handback row for $mdName added to Overview, zh-cn, de-de sheets."
